# Repull data, push all data, mean calculation
# Update column F (dSF) values on Sheet1 to reflect the repulled data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "F2"  = -4
    "F6"  = -3
    "F7"  = 0
    "F8"  = 2
    "F10" = 4
    "F11" = -1
    "F14" = -3
    "F15" = -5
    "F16" = -7
    "F18" = -3
    "F19" = 2
    "F21" = -6
    "F22" = -1
    "F24" = 2
    "F28" = 0
    "F29" = -6
    "F31" = 1
    "F32" = -4
    "F34" = 4
    "F35" = 3
    "F36" = -6
    "F37" = 1
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
